$wb = $excel.ActiveWorkbook

# Add the new "DOS" row to the "Security Attacks" sheet.
# Write column B before column A so the new shared-string entries land in
# the same order as the target workbook (index 10 = "Denial of service
# attack", index 11 = "DOS").
$wsSecurity = $wb.Worksheets.Item("Security Attacks")
$wsSecurity.Range("B3").Value = "Denial of service attack"
$wsSecurity.Range("A3").Value = "DOS"

# Make "Security Attacks" the active sheet/tab again (moves tabSelected +
# the workbook's activeTab away from "Authenications"), and leave the
# selection on the newly added cell B3.
$wsSecurity.Activate()
$wsSecurity.Range("B3").Select()
